$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for B:E columns (SCPC, KSPC, EDC, TOTAL_BCQ) for rows 2..25
$data = @(
    @(25000, 20000, 12000, 57000),
    @(25000, 20000, 0,     45000),
    @(25000, 20000, 0,     45000),
    @(25000, 20000, 0,     45000),
    @(25000, 20000, 0,     45000),
    @(25000, 20000, 0,     45000),
    @(12500, 10000, 0,     22500),
    @(25000, 20000, 0,     45000),
    @(25000, 20000, 0,     45000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000),
    @(25000, 20000, 20000, 65000)
)

# Delete the PEDC column (old column E); this shifts old column F (TOTAL_BCQ) into E
$ws.Range("E1").EntireColumn.Delete()

# Update header for new column E
$ws.Range("E1").Value = "TOTAL_BCQ"

# Write the updated values for rows 2-25, columns B-E
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
}
